$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing scenario values (B2, B3)
$ws.Range("B2").Value = 20
$ws.Range("B3").Value = 50

# Copy the formatting of row 3 (A3:F3) down onto the new rows 4-7
$ws.Range("A3:F3").Copy()
$ws.Range("A4:F7").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new "Objective Function" column (G) for all data rows
$ws.Range("G1").Value = "Objective Function"
$ws.Range("G2").Value = "Maximize Points"
$ws.Range("G3").Value = "Maximize Ownership"
$ws.Range("G4").Value = "Minimize Ownership"
$ws.Range("G5").Value = "Maximize Points"
$ws.Range("G6").Value = "Maximize Ownership"
$ws.Range("G7").Value = "Minimize Ownership"

# New column width for G (closest achievable value to the target 26.6640625)
$ws.Columns.Item(7).ColumnWidth = 25.83

# Fill in the new scenario rows (4-7)
$ws.Range("A4").Value = "Scenario3"
$ws.Range("B4").Value = 50
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "Forced Unconstrained"
$ws.Range("F4").Value = "None"

$ws.Range("A5").Value = "Scenario4"
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "Forced Unconstrained"
$ws.Range("F5").Value = 3

$ws.Range("A6").Value = "Scenario5"
$ws.Range("B6").Value = 50
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = "QB Needs WR"
$ws.Range("F6").Value = "None"

$ws.Range("A7").Value = "Scenario6"
$ws.Range("B7").Value = 50
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "QB Needs WR"
$ws.Range("F7").Value = "None"

# Update selection to match the committed state
$ws.Range("C8").Select()
